# Generate Report for Handoff
# The localization status report is regenerated: the row that used to
# represent "0822cacf-...md" (handed back, in sync with en-US) and the
# row that used to represent "a19ccab5-...md" (ignored dependency) swap
# places, the former now being "Ready for handoff" with fresh
# handoff/handback timestamps for both locales.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = "a19ccab5-e9ca-4cab-ad56-f8ee072f1184.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("A3").Value = "0822cacf-b845-43cd-b6ac-8d79fdd175df.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = "a19ccab5-e9ca-4cab-ad56-f8ee072f1184.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "a19ccab5-e9ca-4cab-ad56-f8ee072f1184.77927fa78aecdc12cbc3d27452998e4801193aa1.zh-cn.xlf"
$ws.Range("D2").Value = "2016-03-08 12:32:42"
$ws.Range("E2").Value = "a19ccab5-e9ca-4cab-ad56-f8ee072f1184.md"
$ws.Range("F2").Value = "a19ccab5-e9ca-4cab-ad56-f8ee072f1184.77927fa78aecdc12cbc3d27452998e4801193aa1.zh-cn.xlf"

$ws.Range("A3").Value = "0822cacf-b845-43cd-b6ac-8d79fdd175df.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "0822cacf-b845-43cd-b6ac-8d79fdd175df.c77dcca1746c842944ac1bea1dd5679f0a71b7d6.zh-cn.xlf"
$ws.Range("D3").Value = "2016-03-08 12:34:08"
$ws.Range("E3").Value = "0822cacf-b845-43cd-b6ac-8d79fdd175df.md"
$ws.Range("F3").Value = "0822cacf-b845-43cd-b6ac-8d79fdd175df.c77dcca1746c842944ac1bea1dd5679f0a71b7d6.zh-cn.xlf"

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = "a19ccab5-e9ca-4cab-ad56-f8ee072f1184.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "a19ccab5-e9ca-4cab-ad56-f8ee072f1184.77927fa78aecdc12cbc3d27452998e4801193aa1.de-de.xlf"
$ws.Range("D2").Value = "2016-03-08 12:32:55"
$ws.Range("E2").Value = "a19ccab5-e9ca-4cab-ad56-f8ee072f1184.md"
$ws.Range("F2").Value = "a19ccab5-e9ca-4cab-ad56-f8ee072f1184.77927fa78aecdc12cbc3d27452998e4801193aa1.de-de.xlf"

$ws.Range("A3").Value = "0822cacf-b845-43cd-b6ac-8d79fdd175df.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "0822cacf-b845-43cd-b6ac-8d79fdd175df.c77dcca1746c842944ac1bea1dd5679f0a71b7d6.de-de.xlf"
$ws.Range("D3").Value = "2016-03-08 12:34:12"
$ws.Range("E3").Value = "0822cacf-b845-43cd-b6ac-8d79fdd175df.md"
$ws.Range("F3").Value = "0822cacf-b845-43cd-b6ac-8d79fdd175df.c77dcca1746c842944ac1bea1dd5679f0a71b7d6.de-de.xlf"
